$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 201.098592
$ws.Range("N2").Value = 603.295776
$ws.Range("O2").Value = 0.7918622805845071
$ws.Range("P2").Value = 0.791862280584507
$ws.Range("Q2").Value = 249.398250727968
$ws.Range("R2").Value = 2244.584256551712
$ws.Range("S2").Value = 0.7918622805845071
$ws.Range("T2").Value = 0.791862280584507

# Row 3
$ws.Range("O3").Value = 0.1414593902976603
$ws.Range("P3").Value = 0.1414593902976603
$ws.Range("S3").Value = 0.1414593902976603
$ws.Range("T3").Value = 0.1414593902976603

# Row 4
$ws.Range("O4").Value = 0.0666783291178327
$ws.Range("P4").Value = 0.06667832911783268
$ws.Range("S4").Value = 0.0666783291178327
$ws.Range("T4").Value = 0.06667832911783268
